$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.069.83'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '1.855.86'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.86'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.23%  '
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0696'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0987'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '2.122.29'
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.42'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.682'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.854.74'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').Value = '35.030.55'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  +25.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.11%  '
$ws.Range('E29').Value = '  +2.78%  '
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.04'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +23.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.781'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.32%  '
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('E38').Value = '  +12.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '92.03'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('E40').Value = '  +5.95%  '
$ws.Range('D41').Value = '1.349.87'
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +65.85%  '
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0540'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.17%  '
$ws.Range('D49').Value = '2.034.98'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +16.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0682'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.66%  '
